$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "VTC.LTR" row (row 16, the Spatial-smoothing-kernel-size
#    row before the preceding rows shift down). Deleting it shifts every row
#    below it up by one, which automatically:
#       - shrinks the table (ListObject) from A1:E22 to A1:E21
#       - shrinks the sheet's used range from A1:E32 to A1:E31
#       - moves old rows 17-32 into rows 16-31
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Delete()

# ---------------------------------------------------------------------------
# 2. Rewrite the Preprocessing block (rows 13-15): the old "Perform linear
#    trend removal" row is gone; the old single "VTC.THP" field is now split
#    into VTC.THP_FMR (row 13) and VTC.THP_VTC (row 14), and the spatial
#    smoothing row moves up to row 15.
# ---------------------------------------------------------------------------

# Row 13: VTC.THP_FMR
$ws.Range("B13").Value = "Temporal high pass filter number of cycles (if included in fmr preprocessing)"
$ws.Range("C13").Value = 2
$d13Plain = "Leave blank or set 0 if THP should not have been performed on FMR. Otherwise, enter the number of cycles that should have been used.`n"
$d13Bold  = "(TYPICAL VALUES: 2 is the BV default or 0 if performing on VTC instead)"
$ws.Range("D13").Value = $d13Plain + $d13Bold
$d13 = $ws.Range("D13")
$d13.Characters($d13Plain.Length + 1, $d13Bold.Length).Font.Bold = $true
$ws.Range("E13").Value = "VTC.THP_FMR"

# Row 14: VTC.THP_VTC
$ws.Range("B14").Value = "Temporal high pass filter number of cycles (if running after preprocessing on vtc)"
$ws.Range("C14").Value = ""
$d14Plain = "Leave blank or set 0 if THP should not be applied or have been applied to VTC. Otherwise enter the number of cycles.`n"
$d14Bold  = "(TYPICAL VALUES: 3 is the BV default or 0 if performed on FMR instead)"
$ws.Range("D14").Value = $d14Plain + $d14Bold
$d14 = $ws.Range("D14")
$d14.Characters($d14Plain.Length + 1, $d14Bold.Length).Font.Bold = $true
$ws.Range("E14").Value = "VTC.THP_VTC"

# Row 15: VTC.SS (spatial smoothing)
$ws.Range("B15").Value = "Spatial smoothing kernel size in mm"
$ws.Range("C15").Value = 6
$d15Plain = "Leave blank or set 0 to skip spatial smoothing. This will also skip MDM generation.`n"
$d15Bold  = "(TYPICAL VALUES: 6 or 8)"
$ws.Range("D15").Value = $d15Plain + $d15Bold
$d15 = $ws.Range("D15")
$d15.Characters($d15Plain.Length + 1, $d15Bold.Length).Font.Bold = $true
$ws.Range("E15").Value = "VTC.SS"

# Rows 13-15 now wrap onto three lines, so they need to be taller.
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 45

# ---------------------------------------------------------------------------
# 3. Clear the stale example value out of EXCLUDE.PARRUN (now row 19 after
#    the shift) - it no longer ships with a default "1-2" value.
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = ""

# ---------------------------------------------------------------------------
# 4. Flip the default for MTN.OVERWRITE (now row 21 after the shift) from
#    TRUE to FALSE.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = $false

# ---------------------------------------------------------------------------
# 5. Restore the view: scrolled down a bit further and selection moved to
#    the spatial-smoothing value cell.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C15").Select()
